# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.671.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.649.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.705.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0749"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "

$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("E33").Value = "  +2.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.274.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.98%  "

$ws.Range("E35").Value = "  +2.25%  "

$ws.Range("E36").Value = "  +5.86%  "

$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "

$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.782.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.89%  "

$ws.Range("E47").Value = "  +3.03%  "

$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.34%  "

$ws.Range("E51").Value = "  -0.55%  "

